$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be read/written as text so numeric-looking values
# (e.g. "0.999") are preserved as strings, then restore original (default) style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.669.18"
$ws.Range("E2").Value = "  -0.14%  "

$ws.Range("D3").Value = "3.102.16"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "398.28"
$ws.Range("E5").Value = "  +3.48%  "

$ws.Range("D6").Value = "103.37"
$ws.Range("E6").Value = "  -0.06%  "

$ws.Range("D7").Value = "0.537"
$ws.Range("E7").Value = "  -1.56%  "

$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  +0.46%  "

$ws.Range("D10").Value = "37.67"
$ws.Range("E10").Value = "  +1.28%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "0.0857"
$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("D13").Value = "3.568.98"
$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").Value = "18.65"
$ws.Range("E14").Value = "  -0.62%  "

$ws.Range("D15").Value = "7.76"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "1.04"
$ws.Range("E16").Value = "  +5.46%  "

$ws.Range("D17").Value = "3.099.87"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "10.66"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("D19").Value = "51.656.57"
$ws.Range("E19").Value = "  -0.25%  "

$ws.Range("D20").Value = "3.20"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").Value = "70.47"
$ws.Range("E23").Value = "  +0.42%  "

$ws.Range("D24").Value = "265.80"
$ws.Range("E24").Value = "  -1.10%  "

$ws.Range("D25").Value = "3.24"
$ws.Range("E25").Value = "  +2.00%  "

$ws.Range("D26").Value = "7.94"
$ws.Range("E26").Value = "  -6.66%  "

$ws.Range("D27").Value = "27.15"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").Value = "7.25"
$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "0.166"
$ws.Range("E30").Value = "  -4.06%  "

$ws.Range("D31").Value = "0.106"
$ws.Range("E31").Value = "  -2.02%  "

$ws.Range("D32").Value = "10.74"
$ws.Range("E32").Value = "  +3.89%  "

$ws.Range("D33").Value = "0.0494"
$ws.Range("E33").Value = "  +11.20%  "

$ws.Range("D34").Value = "36.54"
$ws.Range("E34").Value = "  +5.90%  "

$ws.Range("D36").Value = "49.97"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").Value = "3.37"
$ws.Range("E38").Value = "  -0.78%  "

$ws.Range("D41").Value = "130.31"
$ws.Range("E41").Value = "  +1.37%  "

$ws.Range("D42").Value = "16.83"
$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("D43").Value = "1.86"
$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("D46").Value = "21.93"
$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("E47").Value = "  -1.68%  "

$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  -2.00%  "

$ws.Range("D49").Value = "2.080.34"
$ws.Range("E49").Value = "  +1.55%  "

$ws.Range("D50").Value = "0.0526"
$ws.Range("E50").Value = "  +34.06%  "

$ws.Range("D51").Value = "0.914"
$ws.Range("E51").Value = "  +10.80%  "

# Swap TheGraph / NEARProtocol rows (39 <-> 40) with refreshed data
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "4.06"
$ws.Range("E39").Value = "  +9.45%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.291"
$ws.Range("E40").Value = "  +1.63%  "

# Restore column D to its original (default) style/format
$ws.Range("D2:D51").Style = "Normal"
